# New weekly price record for "Feria Lagunitas de Puerto Montt - Cilantro"
# is inserted at row 151. This pushes the existing rows 151-187 down to
# 152-188 (the last row, old 187, becomes new 188), and the new row 151
# is populated with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 151, shifting rows 151:187 down to 152:188.
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the new data point.
$ws.Range("A151").Value = 4
$ws.Range("B151").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C151").Value = "Los Lagos"
$ws.Range("D151").Value = 44508
$ws.Range("E151").Value = 10
$ws.Range("F151").Value = 100112040
$ws.Range("G151").Value = "Cilantro"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 150
$ws.Range("K151").Value = 9000
$ws.Range("L151").Value = 9000
$ws.Range("M151").Value = 9000
$ws.Range("N151").Value = "$/caja 36 atados"
$ws.Range("O151").Value = "Región Metropolitana"
$ws.Range("P151").Value = 250
$ws.Range("Q151").Value = 36
$ws.Range("R151").Value = "Hortaliza"
